$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (the three IDM CASSEROLE rows), leaving the
# INDOMARET CARD 2020 row (originally row 6) to shift up to row 3.
$ws.Range("A3:F5").EntireRow.Delete()
